$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-format the existing parameter table (B:F, rows 2-11) -----------------
# Previously these cells used the default "General" number format (a couple of
# them, F6/F7, used the built-in 0.00E+00 format). They should now all use the
# built-in "0.00" number format, still centered.
$dataRange = $ws.Range("B2:F11")
$dataRange.HorizontalAlignment = -4108
$dataRange.NumberFormat = "0.00"

# --- Add the new pAUC / CI rows ----------------------------------------------
$ws.Range("A13").Value = "pAUC"
$ws.Range("A14").Value = "pAUC lower 95% CL"
$ws.Range("A15").Value = "pAUC upper 95% CL"

$labelRange = $ws.Range("A13:A15")
$labelRange.HorizontalAlignment = -4108
$labelRange.Font.Bold = $true

$ws.Range("B13").Value = 0.021746073272763299
$ws.Range("B14").Value = 0.0187126054720755
$ws.Range("B15").Value = 0.024982594572062901

$valueRange = $ws.Range("B13:B15")
$valueRange.HorizontalAlignment = -4108
$valueRange.NumberFormat = "0.0000"

# --- Update the selection to match the new active cell -----------------------
$ws.Range("B14").Select()
